# Applies a weekly refresh of the "Alcachofa" price data: the data rows
# (columns D, H, J, K, L, M, N, O, P, Q) get rotated among rows 2,3,4,5,6,9
# following the cycle 2 -> 5 -> 9 -> 4 -> 3 -> 6 -> 2 (i.e. each row ends up
# holding the values that used to belong to the row that follows it in the
# cycle below). Rows 7 and 8 are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the per-record data we need to move around.
$cols = @("D", "H", "J", "K", "L", "M", "N", "O", "P", "Q")

# Snapshot the current ("before") values for every row involved in the cycle.
$rows = @(2, 3, 4, 5, 6, 9)
$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2()
    }
    $snapshot[$r] = $rowData
}

# Destination row -> source row (i.e. destination row receives the values
# that source row used to hold).
$mapping = @{
    2 = 6
    3 = 4
    4 = 9
    5 = 2
    6 = 3
    9 = 5
}

foreach ($dest in $mapping.Keys) {
    $src = $mapping[$dest]
    $data = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$dest").Value = $data[$c]
    }
}
